$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete the extra (trailing) passenger rows (21..35 -> sheet rows 22..36) ---
$ws.Rows("22:36").Delete()

# --- Replace full passenger names with their trimmed name-fragment, and fix up Age ---
# Columns: A=ID, B=Name, C=Sex, D=Age
$names = @("Harris","Bradley","Laina","Heath","Henry","James","Timothy","Leonard","Oscar","Nicholas","Rut","Elizabeth","William","Johan","Adolfina","Hewlett","Eugene","Charles","Julius","Fatima")
$sexes = @("male","female","female","female","male","male","male","male","female","female","female","female","male","male","female","female","male","male","female","female")
$ages  = @(22,38,26,35,35,12,54,2,27,14,4,58,20,39,14,55,2,22,31,22)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = $sexes[$i]
    $ws.Cells.Item($r, 4).Value = $ages[$i]
}

# --- Worksheet view / layout tidy-up to match the uploaded version ---
$ws.Columns.Item(1).ColumnWidth = 8.43
$ws.Columns.Item(2).ColumnWidth = 8.43
$excel.ActiveWindow.Zoom = 100
$ws.Range("A1").Select() | Out-Null
